$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: query terms shifted up, "fees" -> "fee", "free%education" added ---
$ws.Range("D2").Value = "fee"
$ws.Range("D4").Value = "SRC"
$ws.Range("D5").Value = "registration"
$ws.Range("D6").Value = "free%education"
$ws.Range("D7").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("D9").ClearContents()

# --- Column A: new "University Key Words" entries appended ---
$ws.Range("A10").Value = "education"
$ws.Range("A11").Value = "teacher"
$ws.Range("A12").Value = "lecturer"
$ws.Range("A13").Value = "accommodation"

# --- Column C: "gather" renamed to "gathering", new "mass%meeting" action appended ---
$ws.Range("C23").Value = "gathering"
$ws.Range("C29").Value = "mass%meeting"

# --- Selection, matching the saved view state ---
$ws.Range("D11").Select()
